$d = $word.ActiveDocument

# The project timeline table is the only table in the document.
$t = $d.Tables.Item(1)

# Add a new row at the end of the table. Word clones the formatting
# (cell widths, borders, shading) from the row that was previously last.
$newRow = $t.Rows.Add()

$cells = $newRow.Cells
$cells.Item(1).Range.Text = "04/02/2022"
$cells.Item(2).Range.Text = "1 Hour 45 Minutes"
$cells.Item(3).Range.Text = "Simulation Design"
$cells.Item(4).Range.Text = "Created new designs for the functionality of the simulation actions system – including constructs like technology types, how the economy will work etc. etc."

# The "Development Segment" cell for this entry uses a different shading
# colour than the row it was cloned from, so fix it up explicitly.
# (COM colour values are 0xBBGGRR, i.e. byte-reversed RGB -- 99D2F2 -> F2D299)
$cells.Item(3).Shading.BackgroundPatternColor = 0xF2D299

# Two new blank paragraphs are inserted immediately after the table
# (before the existing blank paragraphs that lead up to "ADD MEETINGS").
$tableEnd = $t.Range.End
$afterTable = $d.Range($tableEnd, $tableEnd)
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$afterTable.InsertXML("<w:p $wNs/><w:p $wNs/>")

Write-Output "done"
